$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mapping of row -> [C (DM_Stat), D (P_Value)] new values
$updates = @{
    2  = @(-1.320381201085663,  0.2002806244740265)
    3  = @(-1.06530972133658,   0.2982885416684207)
    4  = @(-1.740115359088417,  0.09581090626266353)
    5  = @(-0.1946396896434873, 0.8474612439213363)
    6  = @(0.5084019544517114,  0.6162318480661122)
    7  = @(0.0373147151087913,  0.970570602629377)
    8  = @(0.9436313984451751,  0.3556064935840384)
    9  = @(-0.6414055482783612, 0.5278856723302827)
    10 = @(0.773720528572617,   0.4473310495342027)
    11 = @(1.542280599425229,   0.1372680417437044)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
}
